$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths
$ws.Columns.Item(2).ColumnWidth = 12.6328125
$ws.Columns.Item(6).ColumnWidth = 13.90625

# Enter the brainstorming notes in the same order the original author typed
# them, so the shared-string table ends up built in the same sequence.
$ws.Range("B4").Value = "Magic Mirror"
$ws.Range("F4").Value = "Phone Control"
$ws.Range("F6").Value = "Heizung  Licht"
$ws.Range("J4").Value = "automatisch"
$ws.Range("J5").Value = "tür licht an "
$ws.Range("J6").Value = "licht aus"
$ws.Range("B6").Value = "alexa"
$ws.Range("F7").Value = "Tisch hoch runter fahren"
$ws.Range("B7").Value = "kalender"
$ws.Range("B8").Value = "wetter"
$ws.Range("B9").Value = "news"
$ws.Range("C8").Value = "mit eigenen mess daten"
$ws.Range("F8").Value = "farbiges licht"
$ws.Range("B10").Value = "touch face recognicion"
$ws.Range("B11").Value = "raspberry pi 4"
$ws.Range("F9").Value = "(iPad?)"
$ws.Range("F10").Value = "rolläden"

# Selection matching the saved sheet view
$ws.Range("H11").Select()
